# Update the "Instructions" worksheet to reflect the new submission
# instructions layout (adds two extra guidance rows, renumbers the
# remaining rows, and bumps the version string).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# The sheet is protected; temporarily unprotect so we can edit cells,
# then restore protection (with the original settings) at the end.
$ws.Unprotect()

# Bump the version number shown on the instructions sheet.
$ws.Range("A2").Value = "Version 1.2.2"

# Make room for two new guidance rows (new rows 5 and 6) by inserting
# before the current row 5. This shifts the old row 5 ("Antibody name"
# header block, previously starting at row 6) down to start at row 8,
# preserving the existing blank-row gap pattern.
$ws.Range("A5:B6").EntireRow.Insert()

# Row 4 text is split: the "Do not edit the other sheets." sentence
# moves out into its own new row (row 6), and a brand new row (row 5)
# is added about not changing the Antibodies header row.
$ws.Range("A4").Value = "Add your antibodies to the 'Antibodies' sheet."
$ws.Range("A5").Value = "Do not change the headers of the 'Antibodies' sheet."
$ws.Range("A6").Value = "Do not edit the other sheets."

# Restore the worksheet protection that was present before the edit.
$ws.Protect()

Write-Host "Instructions sheet updated"
